# "combined model into merge"
#
# The coverage-comparison sheet gains two new weighted "combined model"
# rows (an 0.8*OD+0.2*SS blend and a 0.5/0.5 blend) in the top table, the
# existing "Combined model 0" row is relabelled to the 0.5/0.5 wording, and
# a percentage-formatted mirror of the new 0.8/0.2 row is appended at the
# bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New row 10 ("Combinded model W  (0.8OD+0.2SS)") - inserted right below
# the existing "all_test" row, formatted the same as row 9 above it.
# ---------------------------------------------------------------------
$ws.Range("A9:J9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A10").Value = "Combinded model W  (0.8OD+0.2SS)"
$ws.Range("B10").Value = 0.8
$ws.Range("C10").Value = 22
$ws.Range("D10").Value = 0.90939863840738933
$ws.Range("E10").Value = 0.84080452102829617
$ws.Range("F10").Value = 0.91281748753302627
$ws.Range("G10").Value = 0.092429765512481096
$ws.Range("H10").Value = 0.87533238821645398
$ws.Range("I10").Value = 0.96001247703654569

# ---------------------------------------------------------------------
# Row 9 is relabelled from "Combined model 0" to the new 0.5/0.5 wording
# (written after row 10 so the shared-string table orders the two new
# strings the same way the source workbook does).
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Combined model W (0.5-0.5)"

# ---------------------------------------------------------------------
# New summary row 20 at the bottom of the sheet - same numbers as row 10,
# shown as percentages on a light grey background.
# ---------------------------------------------------------------------
$ws.Range("A20").Value = 0.8
$ws.Range("B20").Value = 22
$ws.Range("C20").Value = 0.90939863840738933
$ws.Range("D20").Value = 0.84080452102829617
$ws.Range("E20").Value = 0.91281748753302627
$ws.Range("F20").Value = 0.092429765512481096
$ws.Range("G20").Value = 0.87533238821645432
$ws.Range("H20").Value = 0.96001247703654569

$ws.Range("C20:D20").Style = "Percent"
$ws.Range("E20").NumberFormat = "0.0%"
$ws.Range("F20").NumberFormat = "0.0%"
$ws.Range("G20").NumberFormat = "0.000%"
$ws.Range("H20").NumberFormat = "0.0%"

$ws.Range("A20").Interior.Color = 15658734
$ws.Range("B20").Interior.Color = 15658734
$ws.Range("C20").Interior.Color = 15658734
$ws.Range("D20").Interior.Color = 15658734
$ws.Range("E20").Interior.Color = 15658734
$ws.Range("F20").Interior.Color = 15658734
$ws.Range("G20").Interior.Color = 15658734
$ws.Range("H20").Interior.Color = 15658734

# ---------------------------------------------------------------------
# Cosmetic: the source edit also left the active selection on G12.
# ---------------------------------------------------------------------
$ws.Range("G12").Select() | Out-Null
